# Modulos de actividades.xlsx - apply edits described by the commit diff.
#
# Summary of the change:
#  - "back-end" sheet: scroll position moved (topLeftCell A13 -> A10).
#  - "ES" sheet: active selection moved (G17 -> I15); a bunch of per-task
#    status cells that used to hold "pm"/"ep"/"u" plus a secondary owner
#    label ("gabriel"/"folege"/...) were all marked as done ("f") and the
#    now-redundant owner/secondary cells were cleared; and three summary
#    cells (B21, D21, F21) got a "% done" formula
#    (COUNTIF(..,"F")/(COUNTA(..)+COUNTBLANK(..))), with D21/F21 getting
#    the same "U" conditional-formatting rule that B21 already had.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "back-end" sheet: just a scroll/viewport change.
# ---------------------------------------------------------------------
$wsBack = $wb.Worksheets.Item("back-end")
$wsBack.Application.Goto($wsBack.Range("A10"), $false)

# ---------------------------------------------------------------------
# 2. "ES" sheet: cell content + selection + conditional formatting.
# ---------------------------------------------------------------------
$wsES = $wb.Worksheets.Item("ES")

# Row 3
$wsES.Range("F3").Value2 = "f"
$wsES.Range("G3").Value2 = ""

# Row 7
$wsES.Range("B7").Value2 = "f"
$wsES.Range("C7").Value2 = ""
$wsES.Range("D7").Value2 = "f"
$wsES.Range("E7").Value2 = ""

# Row 11
$wsES.Range("B11").Value2 = "f"
$wsES.Range("C11").Value2 = ""

# Row 13
$wsES.Range("B13").Value2 = "f"
$wsES.Range("C13").Value2 = ""
$wsES.Range("D13").Value2 = "f"
$wsES.Range("E13").Value2 = ""

# Row 16
$wsES.Range("B16").Value2 = "f"
$wsES.Range("C16").Value2 = ""
$wsES.Range("F16").Value2 = "f"
$wsES.Range("G16").Value2 = ""

# Row 18
$wsES.Range("D18").Value2 = "f"
$wsES.Range("E18").Value2 = ""
$wsES.Range("F18").Value2 = "f"
$wsES.Range("G18").Value2 = ""

# Row 19
$wsES.Range("B19").Value2 = "f"
$wsES.Range("C19").Value2 = ""
$wsES.Range("F19").Value2 = "f"
$wsES.Range("G19").Value2 = ""

# Row 21: "% done" summary formulas for B, D and F columns.
$wsES.Range("B21").Formula = '=(COUNTIF(B2:B20,"F")/(COUNTA(B2:B20)+COUNTBLANK(B2:B20)))'

$wsES.Range("D21").Formula = '=(COUNTIF(D2:D20,"F")/(COUNTA(D2:D20)+COUNTBLANK(D2:D20)))'
$wsES.Range("D21").Style = $wsES.Range("B21").Style
$wsES.Range("D21").HorizontalAlignment = $wsES.Range("B21").HorizontalAlignment
$wsES.Range("D21").NumberFormat = $wsES.Range("B21").NumberFormat

$wsES.Range("F21").Formula = '=(COUNTIF(F2:F20,"F")/(COUNTA(F2:F20)+COUNTBLANK(F2:F20)))'
$wsES.Range("F21").Style = $wsES.Range("B21").Style
$wsES.Range("F21").HorizontalAlignment = $wsES.Range("B21").HorizontalAlignment
$wsES.Range("F21").NumberFormat = $wsES.Range("B21").NumberFormat

# Give D21 / F21 the same "highlight if U" conditional formatting that B21
# already has.
$fcD = $wsES.Range("D21").FormatConditions.Add(1, 3, '="U"')
$fcD.Interior.Color = $wsES.Range("B21").FormatConditions.Item(1).Interior.Color
$fcF = $wsES.Range("F21").FormatConditions.Add(1, 3, '="U"')
$fcF.Interior.Color = $wsES.Range("B21").FormatConditions.Item(1).Interior.Color

# Selection moved from G17 to I15.
$wsES.Activate()
$wsES.Range("I15").Select()

$wb.Save()
